# The workbook contains a weekly price table for Espinaca (Agrícola del Norte
# S.A. de Arica). A new weekly record is inserted at row 11 (right after the
# header + the first 9 existing data rows), pushing every subsequent record
# down by one row. The sheet's used range therefore grows from A1:R71 to
# A1:R72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 11; everything below (old rows 11..71)
# shifts down to rows 12..72, and the sheet dimension is updated automatically.
$ws.Rows(11).Insert()

# Populate the newly inserted row 11 with the new weekly observation.
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value = "Arica y Parinacota"
$ws.Range("D11").Value = 44749
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 100112012
$ws.Range("G11").Value = "Espinaca"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 2000
$ws.Range("L11").Value = 2500
$ws.Range("M11").Value = 2250
$ws.Range("N11").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 750
$ws.Range("Q11").Value = 3
$ws.Range("R11").Value = "Hortaliza"
